$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 125 (shifts existing rows 125..178 down to 126..179)
$ws.Rows("125:125").Insert()

# Populate the newly inserted row 125 with the new weekly data point
$ws.Range("A125").Value = 11
$ws.Range("B125").Value = "Vega Monumental Concepción"
$ws.Range("C125").Value = "Bíobío"
$ws.Range("D125").Value = 45141
$ws.Range("E125").Value = 8
$ws.Range("F125").Value = 100112001
$ws.Range("G125").Value = "Berenjena"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 100
$ws.Range("K125").Value = 8500
$ws.Range("L125").Value = 9000
$ws.Range("M125").Value = 8750
$ws.Range("N125").Value = "$/caja 60 unidades"
$ws.Range("O125").Value = "Región de Arica y Parinacota"
$ws.Range("P125").Value = 146
$ws.Range("Q125").Value = 60
$ws.Range("R125").Value = "Hortaliza"
